# Auto-generated edit script: updates the cryptos price table
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / link / coin-name cells: assign directly ---
$ws.Range("D2").Value = "30.579.88"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.924.68"
$ws.Range("E3").Value = "  +3.56%  "
$ws.Range("E5").Value = "  +4.86%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("E9").Value = "  +6.04%  "
$ws.Range("E10").Value = "  +9.06%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "1.913.91"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("E14").Value = "  +6.61%  "
$ws.Range("E15").Value = "  +5.07%  "
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "30.616.26"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("E18").Value = "  +3.26%  "
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "2.162.65"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("E22").Value = "  +8.61%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +4.71%  "
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("E27").Value = "  +7.74%  "
$ws.Range("E28").Value = "  +11.66%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("E31").Value = "  +4.88%  "
$ws.Range("E32").Value = "  +10.20%  "
$ws.Range("E33").Value = "  +3.23%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E34").Value = "  +4.80%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E35").Value = "  +2.31%  "
$ws.Range("E36").Value = "  +7.98%  "
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("E39").Value = "  +5.25%  "
$ws.Range("E40").Value = "  +5.99%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E42").Value = "  +6.50%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E43").Value = "  +7.69%  "
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("E47").Value = "  +4.97%  "
$ws.Range("E48").Value = "  +17.98%  "
$ws.Range("E49").Value = "  +3.93%  "
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("E51").Value = "  +8.44%  "

# --- Numeric-looking price cells: these must stay TEXT (the sheet
# stores prices as inline strings, e.g. "247.48", not numbers) so we
# flip the cell to text format, assign the literal string, then
# restore the original ("Normal") style so no formatting drifts. ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4738"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06789"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "105.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07733"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.314"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6732"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "288.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007634"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.444"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.330"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.402"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.155"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1083"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.363"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.227"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.191"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05068"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.168"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7436"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02076"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.747"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.691"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.068"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "111.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8844"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.963"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4381"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.282"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.362"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1236"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4051"
$ws.Range("D51").Style = "Normal"
